$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.08528554791034348
    "C2" = 0.4252855479103435
    "D2" = 0.7552855479103435
    "E2" = 0.6731635053256353
    "F2" = 0.3446377147317632
    "G2" = 0.6452855479103434
    "H2" = 0.4212855479103435

    "B3" = 0.34
    "C3" = 0.6699999999999999
    "D3" = 0.5878779574152918
    "E3" = 0.2593521668214197
    "F3" = 0.5600000000000001
    "G3" = 0.336

    "B4" = 0.33
    "C4" = 0.2478779574152918
    "D4" = -0.08064783317858029
    "E4" = 0.22
    "F4" = -0.003999999999999997
    "G4" = 0.1225100404635037
    "H4" = -0.2182525219575302
    "I4" = 0.08717552522494373
    "J4" = -0.1742297805489477

    "B5" = -0.08212204258470818
    "C5" = -0.4106478331785803
    "D5" = -0.11
    "E5" = -0.334
    "F5" = -0.2074899595364962
    "G5" = -0.5482525219575302
    "H5" = -0.2428244747750563
    "I5" = -0.5042297805489477

    "B6" = -0.3285257905938721
    "C6" = -0.02787795741529181
    "D6" = -0.2518779574152918
    "E6" = -0.1253679169517881
    "F6" = -0.466130479372822
    "G6" = -0.1607024321903481
    "H6" = -0.4221077379642395

    "B7" = 0.3006478331785803
    "C7" = 0.0766478331785803
    "D7" = 0.203157873642084
    "E7" = -0.1376046887789499
    "F7" = 0.167823358403524
    "G7" = -0.0935819473703674

    "B8" = -0.224
    "C8" = -0.09748995953649625
    "D8" = -0.4382525219575302
    "E8" = -0.1328244747750563
    "F8" = -0.3942297805489477
    "G8" = -0.4194371574146135
    "H8" = -0.2271788341830432
    "I8" = -0.2759495356205764

    "B9" = 0.1265100404635037
    "C9" = -0.2142525219575302
    "D9" = 0.09117552522494374
    "E9" = -0.1702297805489477
    "F9" = -0.1954371574146135
    "G9" = -0.003178834183043253
    "H9" = -0.0519495356205764

    "B10" = -0.3407625624210339
    "C10" = -0.03533451523856001
    "D10" = -0.2967398210124514
    "E10" = -0.3219471978781172
    "F10" = -0.129688874646547
    "G10" = -0.1784595760840801

    "B11" = 0.3054280471824739
    "C11" = 0.04402274140858248
    "D11" = 0.01881536454291668
    "E11" = 0.2110736877744869
    "F11" = 0.1623029863369538

    "B12" = -0.2614053057738914
    "C12" = -0.2866126826395572
    "D12" = -0.09435435940798698
    "E12" = -0.1431250608455201

    "B13" = -0.02520737686566579
    "C13" = 0.1670509463659045
    "D13" = 0.1182802449283713

    "B14" = 0.1922583232315702
    "C14" = 0.1434876217940371

    "B15" = -0.04877070143753315
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
